$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.296.81"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.933.55"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "357.98"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "110.42"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "19.60"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "3.397.08"
$ws.Range("D16").Value = "2.945.09"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "52.296.90"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "3.54"
$ws.Range("E19").Value = "  +6.72%  "
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "70.59"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "269.42"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "0.183"
$ws.Range("E26").Value = "  +5.28%  "
$ws.Range("E27").Value = "  +15.94%  "
$ws.Range("D28").Value = "27.06"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E30").Value = "  +6.32%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "37.71"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "6.23"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "52.30"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("D39").Value = "18.30"
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("D41").Value = "2.77"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").Value = "23.01"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "119.48"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "3.47"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("D48").Value = "2.134.03"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").Value = "0.929"
$ws.Range("E51").Value = "  -3.91%  "
